# Updated excel with results for par vs seq benchmark
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "par w/ other impl" (P) and "seq w/ other impl" (Q) benchmark results
# for the existing rows 6-15 (rows 6-12 already had A/B/C data, rows 13-15 are new)
$ws.Range("P6").Value  = 69.665000000000006
$ws.Range("Q6").Value  = 7.7560000000000002

$ws.Range("P7").Value  = 63.83
$ws.Range("Q7").Value  = 7.4249999999999998

$ws.Range("P8").Value  = 73.980999999999995
$ws.Range("Q8").Value  = 9.8859999999999992

$ws.Range("P9").Value  = 69.688000000000002
$ws.Range("Q9").Value  = 7.5309999999999997

$ws.Range("P10").Value = 75.709999999999994
$ws.Range("Q10").Value = 7.9059999999999997

$ws.Range("P11").Value = 72.152000000000001
$ws.Range("Q11").Value = 8.5670000000000002

$ws.Range("P12").Value = 69.549000000000007
$ws.Range("Q12").Value = 6.093

$ws.Range("P13").Value = 65.337000000000003
$ws.Range("Q13").Value = 7.4790000000000001

$ws.Range("P14").Value = 64.015000000000001
$ws.Range("Q14").Value = 6.9619999999999997

$ws.Range("P15").Value = 65.867000000000004
$ws.Range("Q15").Value = 6.4580000000000002

# Additional standalone P-only samples in rows 17-26
$ws.Range("P17").Value = 69
$ws.Range("P18").Value = 71
$ws.Range("P19").Value = 67
$ws.Range("P20").Value = 66
$ws.Range("P21").Value = 67
$ws.Range("P22").Value = 63
$ws.Range("P23").Value = 68
$ws.Range("P24").Value = 64
$ws.Range("P25").Value = 66
$ws.Range("P26").Value = 60

# Match the author's final selection/active cell
$ws.Range("Q18").Select()
